$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "Address"

# Re-enter the existing mobile numbers as text (quote-prefixed) so they are
# stored as shared strings instead of numbers
$ws.Range("H2").Value = "'9899001068"
$ws.Range("H3").Value = "'9899001068"
$ws.Range("H4").Value = "'1234567890"

# New row of address data (row 5)
$ws.Range("A5").Value = "Rohan Shukla"
$ws.Range("B5").Value = "Mumbai"
$ws.Range("C5").Value = "mumbai2"
$ws.Range("D5").Value = "Seepz"
$ws.Range("E5").Value = "Andheri"
$ws.Range("F5").Value = "Maharashtra"
$ws.Range("G5").Value = 400096
$ws.Range("H5").Value = "'9876543210"

# Remove the now-unused trailing blank row
$ws.Rows("6").Delete()

# Column A needs to widen slightly to fit "Rohan Shukla"
$ws.Columns("A").ColumnWidth = 12

# Selection moved in the saved file
$ws.Range("E11").Select() | Out-Null
